$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the numeric value in E2 (Codigo Postal)
$ws.Range("E2").Value = 50000

# Update the RFC de la Empresa text in H2
$ws.Range("H2").Value = "UPS7172639"

# Move the active selection to E3, matching the saved sheet view
$ws.Activate()
$ws.Range("E3").Select()
